$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E hold plain text values (prices / percentages).
# Force text number format before assignment so Excel COM does not
# auto-convert the text into a numeric/percentage value, then clear the
# temporary format so the cell keeps its original (default/general) style.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "324.75"
Set-TextValue $ws.Range("E2") "2.44%"
Set-TextValue $ws.Range("D3") "39.85"
Set-TextValue $ws.Range("E3") "5.24%"
Set-TextValue $ws.Range("D4") "5.849"
Set-TextValue $ws.Range("E4") "12.88%"
Set-TextValue $ws.Range("D5") "0.07984"
Set-TextValue $ws.Range("E5") "-0.73%"
Set-TextValue $ws.Range("D6") "4.573"
Set-TextValue $ws.Range("E6") "1.83%"
Set-TextValue $ws.Range("D7") "8.695"
Set-TextValue $ws.Range("E7") "2.12%"
Set-TextValue $ws.Range("D8") "1.916"
Set-TextValue $ws.Range("E8") "-0.39%"
Set-TextValue $ws.Range("D9") "2.941"
Set-TextValue $ws.Range("E9") "-1.95%"
Set-TextValue $ws.Range("D10") "0.9389"
Set-TextValue $ws.Range("E10") "0.01%"
Set-TextValue $ws.Range("D11") "0.1261"
Set-TextValue $ws.Range("E11") "-1.99%"
Set-TextValue $ws.Range("D12") "0.1959"
Set-TextValue $ws.Range("E12") "0.82%"
Set-TextValue $ws.Range("D13") "8.814"
Set-TextValue $ws.Range("E13") "33.79%"
Set-TextValue $ws.Range("D14") "0.09177"
Set-TextValue $ws.Range("E14") "1.02%"
Set-TextValue $ws.Range("D15") "0.03546"
Set-TextValue $ws.Range("E15") "4.34%"
Set-TextValue $ws.Range("D16") "0.09613"
Set-TextValue $ws.Range("E16") "0.72%"
Set-TextValue $ws.Range("D17") "0.001297"
Set-TextValue $ws.Range("E17") "-7.53%"
Set-TextValue $ws.Range("D18") "0.006207"
Set-TextValue $ws.Range("E18") "0.63%"
Set-TextValue $ws.Range("D19") "3.349"
Set-TextValue $ws.Range("E19") "-0.77%"
Set-TextValue $ws.Range("E20") "0.09%"
Set-TextValue $ws.Range("D21") "0.1433"
Set-TextValue $ws.Range("E21") "8.86%"
Set-TextValue $ws.Range("D22") "0.2416"
Set-TextValue $ws.Range("E22") "4.62%"
Set-TextValue $ws.Range("D23") "0.04451"
Set-TextValue $ws.Range("E23") "1.25%"
Set-TextValue $ws.Range("D24") "0.001261"
Set-TextValue $ws.Range("E24") "2.54%"
Set-TextValue $ws.Range("D25") "0.004327"
Set-TextValue $ws.Range("E25") "1.34%"
Set-TextValue $ws.Range("D26") "0.0001144"
Set-TextValue $ws.Range("E26") "-13.91%"
Set-TextValue $ws.Range("D39") "0.02434"
Set-TextValue $ws.Range("E39") "3.43%"
Set-TextValue $ws.Range("D40") "0.05239"
Set-TextValue $ws.Range("E40") "1.48%"
Set-TextValue $ws.Range("D41") "0.007452"
Set-TextValue $ws.Range("E41") "-2.32%"
Set-TextValue $ws.Range("E42") "0.65%"
Set-TextValue $ws.Range("D43") "0.008695"
Set-TextValue $ws.Range("E43") "-0.02%"
Set-TextValue $ws.Range("D44") "0.002128"
Set-TextValue $ws.Range("E44") "0.91%"
Set-TextValue $ws.Range("D45") "0.01052"
Set-TextValue $ws.Range("E45") "26.83%"
Set-TextValue $ws.Range("D46") "0.00006843"
Set-TextValue $ws.Range("E46") "5.81%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.37%"
Set-TextValue $ws.Range("D48") "0.002878"
Set-TextValue $ws.Range("E48") "0.64%"
Set-TextValue $ws.Range("D49") "0.001425"
Set-TextValue $ws.Range("E49") "-15.61%"
Set-TextValue $ws.Range("D50") "0.00002107"
Set-TextValue $ws.Range("E50") "0.37%"
Set-TextValue $ws.Range("D51") "0.0002007"
Set-TextValue $ws.Range("E51") "0.37%"
